$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the stray formatted-but-empty A1 cell (header row, column A) ---
$ws.Range("A1").Clear()

# --- Existing data rows (2-89): drop the now-redundant explicit cell style on
#     columns B/C so they fall back to the sheet's default "Normal" style ---
$ws.Range("B2:C89").Style = "Normal"

# --- Append the 7 new reference rows (index continues 88..94) ---
$newRows = @(
    @(88, "国家或地区的经济政策不确定性(EPU)数据", "article_epu_index"),
    @(89, "碳排放权-国内", "energy_carbon_domestic"),
    @(90, "北京市碳排放权电子交易平台-北京市碳排放权公开交易行情", "energy_carbon_bj"),
    @(91, "碳排放权-深圳", "energy_carbon_sz"),
    @(92, "碳排放权-国际", "energy_carbon_eu"),
    @(93, "碳排放权-湖北", "energy_carbon_hb"),
    @(94, "碳排放权-广州", "energy_carbon_gz")
)

$startRow = 90
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    # Column A keeps the same "index" style used by the rows above it.
    $ws.Range("A89").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Range("A$r").Value = $data[0]

    # Columns B/C use the workbook's default style (same as the rows above
    # after their style cleanup).
    $ws.Range("B$r").Value = $data[1]
    $ws.Range("B$r").Style = "Normal"
    $ws.Range("C$r").Value = $data[2]
    $ws.Range("C$r").Style = "Normal"
}

$excel.CutCopyMode = 0
